$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 2")

function Swap-Rows($Row1, $Row2, $FirstCol, $LastCol) {
    for ($c = $FirstCol; $c -le $LastCol; $c++) {
        $cell1 = $ws.Cells.Item($Row1, $c)
        $cell2 = $ws.Cells.Item($Row2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-Rows 66 67 2 30
Swap-Rows 137 138 2 30
Swap-Rows 155 156 2 30
